# Build site at 2022-08-26 16:06:42 UTC
# Insert a new row for professor "1176388 - Luiz Tadeu Fernandes Eleno"
# right after the existing "519033 - Carlos Yujiro Shigue" row (row 13),
# shifting the "Programa resumido:" block and everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14 (pushes rows 14..24 down to 15..25).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new professor's name. The
# inserted row already inherits the correct per-column formatting (style 2
# for column B, style 3 for column C) from row 13 above it.
$ws.Cells.Item(14, 2).Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Cells.Item(14, 3).Value = "1176388 - Luiz Tadeu Fernandes Eleno"
